$d = $word.ActiveDocument
$wdFindContinue = 1
$wdReplaceOne = 2

$old1 = @'
Statics and dynamics of fluids, oscillations and mechanical waves, ideal gas,temperature, heat and the laws of thermodynamics.
'@
$new1 = @'
Statics and dynamics of fluids, oscillations and mechanical waves, ideal gas,^ltemperature, heat and the laws of thermodynamics.
'@
$found1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new1, $wdReplaceOne)
Write-Output "block1: $found1"

$old2 = @'
1) Estática de fluidos: pressão, princípios de Pascal e Arquimedes, tensão superficial, capilaridade;2) Dinâmica de fluidos: vazão, fluidos ideais, equação da continuidade, equação de Bernoulli, viscosidade, lei de Hagen-Poiseuille;3) Oscilações: movimento harmônico simples, amortecido e forçado, ressonância; 4) Ondas: transversais e longitudinais, equação de onda, superposição, interferência, ondas estacionárias e ressonância, ondas sonoras, intensidade e nível sonoro, batimentos, efeito Doppler;5) Temperatura e calor: conceitos, escalas de temperatura, a lei zero da termodinâmica, dilatação térmica, absorção de calor por sólidos e líquidos, calor e trabalho, mecanismos de transferência de calor, gases ideais, calor específico molar de um gás ideal e graus de liberdade;6) Termodinâmica: primeira lei da termodinâmica, processos reversíveis eirreversíveis, entropia, segunda lei da termodinâmica, máquinas térmicas eeficiência.
'@
$new2 = @'
1) Estática de fluidos: pressão, princípios de Pascal e Arquimedes, tensão superficial, capilaridade;^l2) Dinâmica de fluidos: vazão, fluidos ideais, equação da continuidade, equação de Bernoulli, viscosidade, lei de Hagen-Poiseuille;^l3) Oscilações: movimento harmônico simples, amortecido e forçado, ressonância; ^l4) Ondas: transversais e longitudinais, equação de onda, superposição, interferência, ondas estacionárias e ressonância, ondas sonoras, intensidade e nível sonoro, batimentos, efeito Doppler;^l5) Temperatura e calor: conceitos, escalas de temperatura, a lei zero da termodinâmica, dilatação térmica, absorção de calor por sólidos e líquidos, calor e trabalho, mecanismos de transferência de calor, gases ideais, calor específico molar de um gás ideal e graus de liberdade;^l6) Termodinâmica: primeira lei da termodinâmica, processos reversíveis e^lirreversíveis, entropia, segunda lei da termodinâmica, máquinas térmicas e^leficiência.
'@
$found2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new2, $wdReplaceOne)
Write-Output "block2: $found2"

$old3 = @'
1) Fluids at rest: pressure, Pascal’s principle , Archimedes’ principle, surfacetension and capillarity;2) Fluids in motion: flow rate, ideal fluids, the equation of continuity, Bernoulli’s equation, viscosity and the Hagen-Poiseuille law;3) Oscillation: simple harmonic motion, damped and forced oscillations, resonance;4) Waves: transverse and longitudinal, wave equation, superposition, interference, standing waves, sound waves, intensity and sound level, beats, Doppler effect;5) Temperature and heat: definitions, zeroth Law of thermodynamics, thermal expansion, absorption of heat by solids and liquids, heat and work, heat transfer mechanisms, ideal gases, specific heat and degrees of freedom for an ideal gas;6) Thermodynamics: the first law of thermodynamics, reversible and irreversible processes, heat engines and efficiency, entropy, the second law of thermodynamics.
'@
$new3 = @'
1) Fluids at rest: pressure, Pascal’s principle , Archimedes’ principle, surface^ltension and capillarity;^l2) Fluids in motion: flow rate, ideal fluids, the equation of continuity, Bernoulli’s equation, viscosity and the Hagen-Poiseuille law;^l3) Oscillation: simple harmonic motion, damped and forced oscillations, resonance;^l4) Waves: transverse and longitudinal, wave equation, superposition, interference, standing waves, sound waves, intensity and sound level, beats, Doppler effect;^l5) Temperature and heat: definitions, zeroth Law of thermodynamics, thermal expansion, absorption of heat by solids and liquids, heat and work, heat transfer mechanisms, ideal gases, specific heat and degrees of freedom for an ideal gas;^l6) Thermodynamics: the first law of thermodynamics, reversible and irreversible processes, heat engines and efficiency, entropy, the second law of thermodynamics.
'@
$found3 = $d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new3, $wdReplaceOne)
Write-Output "block3: $found3"

$old4 = @'
NUSSENZVEIG, H.M. Curso de Física Básica. Vol. 2, Edgard Blucher (2008).RESNICK, R.; HALLIDAY, D. Fundamentos de Física. Vol.2, LTC (2008).TIPLER, P.; MOSCA, G. Física para Cientistas e Engenheiros. Vol.2, LTC (2008).SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. Física I, Vol. 2, Pearson Addison Wesley (2009).JEWETT Jr, John W.; SERWAY, Raymond A. Princípios de Física. Vol. 2, Thomson Pioneira (2008).
'@
$new4 = @'
NUSSENZVEIG, H.M. Curso de Física Básica. Vol. 2, Edgard Blucher (2008).^lRESNICK, R.; HALLIDAY, D. Fundamentos de Física. Vol.2, LTC (2008).^lTIPLER, P.; MOSCA, G. Física para Cientistas e Engenheiros. Vol.2, LTC (2008).^lSEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. Física I, Vol. 2, Pearson Addison Wesley (2009).^lJEWETT Jr, John W.; SERWAY, Raymond A. Princípios de Física. Vol. 2, Thomson Pioneira (2008).
'@
$found4 = $d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new4, $wdReplaceOne)
Write-Output "block4: $found4"

if (-not ($found1 -and $found2 -and $found3 -and $found4)) {
    throw "One or more Find/Replace operations failed to locate their target text."
}

